$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, shifting existing rows 16-126 down to 17-127
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with fresh data (matches the dataset pattern)
$ws.Range("A16").Value = 4
$ws.Range("B16").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C16").Value = "Los Lagos"
$ws.Range("D16").Value = 44532
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = 100112009
$ws.Range("G16").Value = "Acelga"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 50
$ws.Range("K16").Value = 3500
$ws.Range("L16").Value = 3500
$ws.Range("M16").Value = 3500
$ws.Range("N16").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O16").Value = "Región del Maule"
$ws.Range("P16").Value = 875
$ws.Range("Q16").Value = 4
$ws.Range("R16").Value = "Hortaliza"
